# Daily attendance processing - normalize "Recorded By" ordering.
# For every "Recorded By" (column G) cell that contains a comma-separated
# list including the exact (case-sensitive) token "System", reverse the
# order of the items in the list. Cells without a "System" token, or with
# only a single value, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -eq $null) {
        continue
    }
    if ($val -eq "") {
        continue
    }

    $parts = $val.Split(",")
    if ($parts.Count -le 1) {
        continue
    }

    $trimmedParts = @()
    for ($j = 0; $j -lt $parts.Count; $j++) {
        $trimmedParts += $parts[$j].Trim()
    }

    $hasSystem = $false
    for ($j = 0; $j -lt $trimmedParts.Count; $j++) {
        if ($trimmedParts[$j].CompareTo("System") -eq 0) {
            $hasSystem = $true
        }
    }

    if ($hasSystem) {
        $reversed = @()
        for ($j = $trimmedParts.Count - 1; $j -ge 0; $j--) {
            $reversed += $trimmedParts[$j]
        }
        $newVal = $reversed -join ", "
        $cell.Value = $newVal
    }
}
